# Update Excel file with latest predictions
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Home win" -- 4 data rows now (was 6)
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Home win")

$data = @(
    @("28-01-2025 23:30","WORLD","SUDAMERICANO U20","Colombia U20 - Ecuador U20",73.3,2.2),
    @("29-01-2025 22:00","BRAZIL","GAÚCHO - 1","Ypiranga-RS - São Luiz",86.7,2.25),
    @("29-01-2025 19:45","ENGLAND","NATIONAL LEAGUE CUP","Braintree - Tottenham Hotspur U21",70,2.25),
    @("29-01-2025 15:00","WORLD","FRIENDLIES CLUBS","St Patrick's Athl. - Vancouver Whitecaps",80,3.6)
)

$row = 2
foreach ($r in $data) {
    $ws.Range("A$row").Value = $r[0]
    $ws.Range("B$row").Value = $r[1]
    $ws.Range("C$row").Value = $r[2]
    $ws.Range("D$row").Value = $r[3]
    $ws.Range("E$row").Value = $r[4]
    $ws.Range("F$row").Value = $r[5]
    $row++
}
# sheet used to have 6 data rows (through row 7); remove the now-unused tail
$ws.Range("A6:F7").Clear()

# ---------------------------------------------------------------------------
# Sheet "Away Win" -- unchanged
# ---------------------------------------------------------------------------

# ---------------------------------------------------------------------------
# Sheet "Draw" -- still 3 data rows, values updated
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Draw")

$data = @(
    @("29-01-2025 22:30","BRAZIL","PAULISTA - A1","Noroeste - Inter De Limeira",80,3),
    @("29-01-2025 22:30","BRAZIL","PAULISTA - A2","Grêmio Prudente - Primavera SP",60,2.95),
    @("29-01-2025 23:00","BRAZIL","PERNAMBUCANO - 1","Sport Recife - Jaguaré",70,6.5)
)

$row = 2
foreach ($r in $data) {
    $ws.Range("A$row").Value = $r[0]
    $ws.Range("B$row").Value = $r[1]
    $ws.Range("C$row").Value = $r[2]
    $ws.Range("D$row").Value = $r[3]
    $ws.Range("E$row").Value = $r[4]
    $ws.Range("F$row").Value = $r[5]
    $row++
}

# ---------------------------------------------------------------------------
# Sheet "Btts" -- 10 data rows now (was 7)
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Btts")

$data = @(
    @("29-01-2025 20:00","WORLD","UEFA CHAMPIONS LEAGUE","Aston Villa - Celtic",80,1.7),
    @("29-01-2025 20:00","WORLD","UEFA CHAMPIONS LEAGUE","Bayern München - Slovan Bratislava",80,3),
    @("29-01-2025 20:00","WORLD","UEFA CHAMPIONS LEAGUE","Inter - Monaco",76.7,1.7),
    @("29-01-2025 20:00","WORLD","UEFA CHAMPIONS LEAGUE","Lille - Feyenoord",76,1.7),
    @("29-01-2025 20:00","WORLD","UEFA CHAMPIONS LEAGUE","Manchester City - Club Brugge KV",83.3,1.95),
    @("29-01-2025 21:30","BRAZIL","PAULISTA - A1","São Bernardo - Santos",76.7,2.05),
    @("29-01-2025 23:00","BRAZIL","PERNAMBUCANO - 1","Afogados - Central SC",80,1.95),
    @("29-01-2025 23:00","BRAZIL","PERNAMBUCANO - 1","Sport Recife - Jaguaré",80,2.5),
    @("29-01-2025 11:00","PORTUGAL","LIGA REVELAÇÃO U23","Farense U23 - Portimonense U23",76.7,1.73),
    @("26-01-2025 11:00","SPAIN","SEGUNDA DIVISIÓN RFEF - GROUP 2","Real Sociedad III - Tudelano",83.3,1.85)
)

$row = 2
foreach ($r in $data) {
    $ws.Range("A$row").Value = $r[0]
    $ws.Range("B$row").Value = $r[1]
    $ws.Range("C$row").Value = $r[2]
    $ws.Range("D$row").Value = $r[3]
    $ws.Range("E$row").Value = $r[4]
    $ws.Range("F$row").Value = $r[5]
    $row++
}

# ---------------------------------------------------------------------------
# Sheet "Over_Under" -- 6 data rows now (was 9)
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Over_Under")

$data = @(
    @("29-01-2025 20:00","WORLD","UEFA CHAMPIONS LEAGUE","Inter - Monaco",70,1.67,65,2.62),
    @("29-01-2025 20:00","WORLD","UEFA CHAMPIONS LEAGUE","VfB Stuttgart - Paris Saint Germain",80,2,53.3,3.2),
    @("29-01-2025 15:00","MEXICO","U23 LEAGUE","Atlas U23 - Monterrey U23",80,1.75,10,2.88),
    @("29-01-2025 15:00","MEXICO","U23 LEAGUE","Queretaro U23 - Pachuca U23",70,1.83,60,3.1),
    @("29-01-2025 17:00","PORTUGAL","LIGA REVELAÇÃO U23","Benfica U23 - Estrela U23",70,1.62,60,2.55),
    @("29-01-2025 10:00","WORLD","FRIENDLIES CLUBS","Hrvace - Radnik Bijeljina",73.3,1.5,66.7,2.25)
)

$row = 2
foreach ($r in $data) {
    $ws.Range("A$row").Value = $r[0]
    $ws.Range("B$row").Value = $r[1]
    $ws.Range("C$row").Value = $r[2]
    $ws.Range("D$row").Value = $r[3]
    $ws.Range("E$row").Value = $r[4]
    $ws.Range("F$row").Value = $r[5]
    $ws.Range("G$row").Value = $r[6]
    $ws.Range("H$row").Value = $r[7]
    $row++
}
# sheet used to have 9 data rows (through row 10); remove the now-unused tail
$ws.Range("A8:H10").Clear()

Write-Host "Done updating predictions"
